$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.221.88"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "1.870.07"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'234.56"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'0.4702"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  -1.42%  "

$ws.Range("D9").Value = "'41.68"
$ws.Range("E9").Value = "  -2.42%  "

$ws.Range("D10").Value = "'0.06551"
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("D11").Value = "'21.20"
$ws.Range("E11").Value = "  -2.21%  "

$ws.Range("D12").Value = "'0.07812"
$ws.Range("E12").Value = "  -1.72%  "

$ws.Range("D13").Value = "'96.48"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").Value = "1.850.07"
$ws.Range("E14").Value = "  -0.03%  "

$ws.Range("D15").Value = "'0.6897"
$ws.Range("E15").Value = "  +2.41%  "

$ws.Range("E16").Value = "  +0.11%  "

$ws.Range("D17").Value = "'266.49"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "30.210.76"

$ws.Range("D19").Value = "'13.74"
$ws.Range("E19").Value = "  +1.20%  "

$ws.Range("D20").Value = "'0.000007704"
$ws.Range("E20").Value = "  +1.86%  "

$ws.Range("D22").Value = "2.102.29"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "'5.233"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("D25").Value = "'6.153"
$ws.Range("E25").Value = "  +0.46%  "

$ws.Range("D26").Value = "'9.499"
$ws.Range("E26").Value = "  +3.91%  "

$ws.Range("D27").Value = "'165.61"
$ws.Range("E27").Value = "  -0.60%  "

$ws.Range("D28").Value = "'18.72"
$ws.Range("E28").Value = "  -0.37%  "

$ws.Range("D29").Value = "'1.932"
$ws.Range("E29").Value = "  +0.43%  "

$ws.Range("E30").Value = "  -1.63%  "

$ws.Range("D31").Value = "'0.09918"

$ws.Range("D32").Value = "'4.344"
$ws.Range("E32").Value = "  +1.86%  "

$ws.Range("D33").Value = "'1.453"
$ws.Range("E33").Value = "  -0.75%  "

$ws.Range("E34").Value = "  +1.20%  "

$ws.Range("D35").Value = "'0.04733"
$ws.Range("E35").Value = "  +1.16%  "

$ws.Range("D36").Value = "'1.128"
$ws.Range("E36").Value = "  +1.28%  "

$ws.Range("D37").Value = "'0.6994"
$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("D38").Value = "'2.718"
$ws.Range("E38").Value = "  +0.33%  "

$ws.Range("D39").Value = "'0.01861"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("D40").Value = "'2.786"
$ws.Range("E40").Value = "  +7.46%  "

$ws.Range("D41").Value = "'6.253"
$ws.Range("E41").Value = "  -1.12%  "

$ws.Range("D42").Value = "'72.54"
$ws.Range("E42").Value = "  -0.68%  "

$ws.Range("D43").Value = "'1.936"
$ws.Range("E43").Value = "  +0.53%  "

$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D46").Value = "'0.8322"
$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("D47").Value = "'102.91"
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("D48").Value = "'977.60"
$ws.Range("E48").Value = "  +4.22%  "

$ws.Range("D49").Value = "'7.075"
$ws.Range("E49").Value = "  +1.39%  "

$ws.Range("D50").Value = "'9.138"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("E51").Value = "  +2.03%  "
